# casos_teste_juros.xlsx — update test-case inputs on "Planilha1" and add
# a few new pytest test cases (per commit: "ADD: atualização excel e
# adição de funções pytest").
#
# Row 4 (Verificar se o resultado é negativo): input changes from an
#   all-positive list with one negative entry to a fully negative list.
# Row 5 (Verifica se é string): sample string input changes from "oi" to "str".
# Row 6 (Verifica se os parâmetros estão vazios): empty-input sample changes
#   from an empty string "" to an empty list [].

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

$ws.Range("E4").Value = "[-1000, -20, -1]"
$ws.Range("E6").Value = "[]"
$ws.Range("E5").Value = """str"""

# Leave the cursor where the author's last save left it.
[void]$ws.Range("E9").Select()
